$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.283.04'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.910.39'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.720'
$ws.Range('E5').Value = '  +9.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '254.13'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.51'
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.365'
$ws.Range('E9').Value = '  +4.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.51'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0762'
$ws.Range('E11').Value = '  +6.20%  '
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.187.81'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.72'
$ws.Range('E14').Value = '  +4.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.719'
$ws.Range('E15').Value = '  +2.23%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.943.23'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.92'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.310.50'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.28'
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0849'
$ws.Range('E20').Value = '  +3.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.57'
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('E22').Value = '  +4.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.10'
$ws.Range('E23').Value = '  +5.10%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  +4.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.38'
$ws.Range('E26').Value = '  +3.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.17'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.69'
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.74'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('E30').Value = '  +4.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.126.75'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.36'
$ws.Range('E32').Value = '  +4.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  +13.99%  '
$ws.Range('E34').Value = '  +23.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0585'
$ws.Range('E35').Value = '  +3.04%  '
$ws.Range('E36').Value = '  +1.96%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.913'
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0217'
$ws.Range('E40').Value = '  +4.42%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.08'
$ws.Range('E41').Value = '  +4.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.60'
$ws.Range('E42').Value = '  +7.15%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0654'
$ws.Range('E44').Value = '  +1.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.338.19'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.41'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('E48').Value = '  +2.75%  '
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.34'
$ws.Range('E50').Value = '  -5.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.97'
$ws.Range('E51').Value = '  +16.14%  '
